# Re-pull / push updated dSF (column F) data for each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = -5
    4  = 0
    5  = 2
    6  = 1
    7  = -1
    8  = -7
    9  = -1
    10 = 5
    11 = 5
    12 = -2
    13 = -1
    14 = 1
    15 = 2
    16 = 0
    17 = 0
    18 = -3
    19 = 0
    20 = -5
    21 = 0
    22 = 2
    23 = 0
    24 = -3
    25 = -4
    26 = 7
    27 = -3
    28 = 3
    29 = 0
    30 = 0
    31 = -5
    32 = -1
    33 = 2
    34 = -3
    35 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
